# Refresh the crypto symbol list: update Price (column D) and
# Volume(1h) (column E) figures for each coin row, matching the
# latest scrape snapshot (GitHub Actions run, 2023-01-07 14:42 UTC).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Cell, $NewValue) {
    $rng = $Sheet.Range($Cell)
    # Force text storage so numeric/percent-looking strings are not
    # reinterpreted as numbers, then drop back to the default cell
    # style so only the value itself changes.
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" '260.35'
Set-TextValue $ws "E2" '1.64%'
Set-TextValue $ws "D3" '27.23'
Set-TextValue $ws "E3" '2.41%'
Set-TextValue $ws "D4" '4.686'
Set-TextValue $ws "E4" '1.03%'
Set-TextValue $ws "D5" '0.06159'
Set-TextValue $ws "E5" '3.94%'
Set-TextValue $ws "D6" '6.672'
Set-TextValue $ws "E6" '1.05%'
Set-TextValue $ws "D7" '0.8507'
Set-TextValue $ws "E7" '-0.60%'
Set-TextValue $ws "D8" '0.9156'
Set-TextValue $ws "E8" '0.50%'
Set-TextValue $ws "E9" '2.13%'
Set-TextValue $ws "D10" '0.04811'
Set-TextValue $ws "E10" '10.12%'
Set-TextValue $ws "D11" '0.07084'
Set-TextValue $ws "E11" '1.13%'
Set-TextValue $ws "D12" '0.03104'
Set-TextValue $ws "E12" '2.83%'
Set-TextValue $ws "E13" '-0.53%'
Set-TextValue $ws "D14" '0.001533'
Set-TextValue $ws "E14" '0.11%'
Set-TextValue $ws "D15" '0.0006123'
Set-TextValue $ws "E15" '1.04%'
Set-TextValue $ws "D16" '0.006038'
Set-TextValue $ws "E16" '0.13%'
Set-TextValue $ws "D17" '3.452'
Set-TextValue $ws "E17" '-0.66%'
Set-TextValue $ws "E18" '0.80%'
Set-TextValue $ws "E19" '1.39%'
Set-TextValue $ws "D20" '0.3078'
Set-TextValue $ws "E20" '-0.05%'
Set-TextValue $ws "D21" '0.1296'
Set-TextValue $ws "E21" '0.85%'
Set-TextValue $ws "D22" '4.088'
Set-TextValue $ws "E22" '5.12%'
Set-TextValue $ws "D23" '0.04219'
Set-TextValue $ws "E23" '0.20%'
Set-TextValue $ws "D24" '0.001216'
Set-TextValue $ws "E24" '-0.09%'
Set-TextValue $ws "D25" '0.003801'
Set-TextValue $ws "E25" '-17.83%'
Set-TextValue $ws "E26" '-0.03%'
Set-TextValue $ws "E27" '-8.02%'
Set-TextValue $ws "D40" '0.03874'
Set-TextValue $ws "E40" '1.99%'
Set-TextValue $ws "D41" '0.1113'
Set-TextValue $ws "E41" '1.15%'
Set-TextValue $ws "D42" '0.004086'
Set-TextValue $ws "E42" '-34.34%'
Set-TextValue $ws "E43" '13.68%'
Set-TextValue $ws "E44" '-4.44%'
Set-TextValue $ws "D45" '0.00005157'
Set-TextValue $ws "E45" '0.73%'
Set-TextValue $ws "E46" '-0.04%'
Set-TextValue $ws "E47" '8.20%'
Set-TextValue $ws "D48" '0.1654'
Set-TextValue $ws "E48" '-31.32%'
Set-TextValue $ws "E49" '-0.04%'
Set-TextValue $ws "E50" '-0.04%'
